# Auto-generated Excel COM-interop script to apply market-price / profit updates
# to the Cactuar_Profits workbook, per the scheduled-runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 516.44446
$ws.Range("I2").Value = 414
$ws.Range("K2").Value = 414
$ws.Range("M2").Value = -301

$ws.Range("H94").Value = 11117382
$ws.Range("I94").Value = 15877118
$ws.Range("K94").Value = 15877118
$ws.Range("M94").Value = -15876667

$ws.Range("H121").Value = 4342.52
$ws.Range("J121").Value = 4342.52
$ws.Range("L121").Value = 13027.56
$ws.Range("N121").Value = -16521.56

$ws.Range("H135").Value = 4451.5386
$ws.Range("I135").Value = 1633.7142
$ws.Range("J135").Value = 7739
$ws.Range("K135").Value = 14703.4278
$ws.Range("L135").Value = 69651
$ws.Range("M135").Value = -12168.4278
$ws.Range("N135").Value = -74721

$ws.Range("H137").Value = 11841702
$ws.Range("I137").Value = 714741.4399999999
$ws.Range("J137").Value = 22226866
$ws.Range("K137").Value = 2144224.32
$ws.Range("L137").Value = 66680598
$ws.Range("M137").Value = -2141674.32
$ws.Range("N137").Value = -66685698

$ws.Range("H138").Value = 5097.5884
$ws.Range("J138").Value = 7098.121
$ws.Range("L138").Value = 21294.363
$ws.Range("N138").Value = -31574.363

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14802.766
$ws.Range("I32").Value = 15587.902
$ws.Range("K32").Value = 15587.902
$ws.Range("M32").Value = -15300.902

$ws.Range("H61").Value = 9781.962
$ws.Range("I61").Value = 12301.789
$ws.Range("K61").Value = 12301.789
$ws.Range("M61").Value = -12089.789

$ws.Range("H110").Value = 1078462.2
$ws.Range("I110").Value = 1702668.2
$ws.Range("K110").Value = 1702668.2
$ws.Range("M110").Value = -1700623.2

$ws.Range("H136").Value = 9781.962
$ws.Range("I136").Value = 12301.789
$ws.Range("K136").Value = 36905.367
$ws.Range("M136").Value = -34355.367

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2105.5518
$ws.Range("I20").Value = 2080.6843
$ws.Range("J20").Value = 2152.8
$ws.Range("K20").Value = 2080.6843
$ws.Range("L20").Value = 2152.8
$ws.Range("M20").Value = -1833.6843
$ws.Range("N20").Value = -2646.8

$ws.Range("H134").Value = 2975.7307
$ws.Range("I134").Value = 1189.2632
$ws.Range("K134").Value = 3567.7896
$ws.Range("M134").Value = -1032.7896

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2255.158
$ws.Range("I16").Value = 2250.0908
$ws.Range("K16").Value = 2250.0908
$ws.Range("M16").Value = -1963.0908

$ws.Range("H22").Value = 1173.4
$ws.Range("I22").Value = 901.4
$ws.Range("J22").Value = 1445.4
$ws.Range("K22").Value = 901.4
$ws.Range("L22").Value = 1445.4
$ws.Range("M22").Value = -551.4
$ws.Range("N22").Value = -2145.4

$ws.Range("H31").Value = 4390.4106
$ws.Range("I31").Value = 841.34784
$ws.Range("J31").Value = 6864
$ws.Range("K31").Value = 841.34784
$ws.Range("L31").Value = 6864
$ws.Range("M31").Value = -546.34784
$ws.Range("N31").Value = -7454

$ws.Range("H34").Value = 4390.4106
$ws.Range("I34").Value = 841.34784
$ws.Range("J34").Value = 6864
$ws.Range("K34").Value = 841.34784
$ws.Range("L34").Value = 6864
$ws.Range("M34").Value = -639.34784
$ws.Range("N34").Value = -7268

$ws.Range("H103").Value = 34583.168
$ws.Range("I103").Value = 12008
$ws.Range("J103").Value = 57158.332
$ws.Range("K103").Value = 12008
$ws.Range("L103").Value = 57158.332
$ws.Range("M103").Value = -10836
$ws.Range("N103").Value = -59502.332

$ws.Range("H113").Value = 2255.158
$ws.Range("I113").Value = 2250.0908
$ws.Range("K113").Value = 2250.0908
$ws.Range("M113").Value = -80.09079999999994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 277.91306
$ws.Range("I2").Value = 170.18182
$ws.Range("J2").Value = 376.66666
$ws.Range("K2").Value = 1021.09092
$ws.Range("L2").Value = 2259.99996
$ws.Range("M2").Value = -908.0909199999999
$ws.Range("N2").Value = -2485.99996

$ws.Range("H5").Value = 855.5714
$ws.Range("I5").Value = 918.4
$ws.Range("K5").Value = 2755.2
$ws.Range("M5").Value = -2643.2

$ws.Range("H37").Value = 77029464
$ws.Range("J37").Value = 77029464
$ws.Range("L37").Value = 231088392
$ws.Range("N37").Value = -231088616

$ws.Range("H38").Value = 236.0625
$ws.Range("J38").Value = 295.9091
$ws.Range("L38").Value = 887.7273
$ws.Range("N38").Value = -1581.7273

$ws.Range("H94").Value = 7930
$ws.Range("J94").Value = 7930
$ws.Range("L94").Value = 23790
$ws.Range("N94").Value = -25142

$ws.Range("H113").Value = 877.2
$ws.Range("I113").Value = 792.3333
$ws.Range("J113").Value = 1004.5
$ws.Range("K113").Value = 2376.9999
$ws.Range("L113").Value = 3013.5
$ws.Range("M113").Value = -206.9998999999998
$ws.Range("N113").Value = -7353.5

$ws.Range("H134").Value = 8439.263000000001
$ws.Range("I134").Value = 2026.7333
$ws.Range("K134").Value = 6080.199900000001
$ws.Range("M134").Value = -1010.199900000001

$ws.Range("H135").Value = 855.5714
$ws.Range("I135").Value = 918.4
$ws.Range("K135").Value = 8265.6
$ws.Range("M135").Value = -5730.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 52858.69
$ws.Range("I132").Value = 78771.92999999999
$ws.Range("J132").Value = 6214.8667
$ws.Range("K132").Value = 236315.79
$ws.Range("L132").Value = 18644.6001
$ws.Range("M132").Value = -233785.79
$ws.Range("N132").Value = -23704.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1600
$ws.Range("I2").Value = 857.1429000000001
$ws.Range("J2").Value = 3333.3333
$ws.Range("K2").Value = 857.1429000000001
$ws.Range("L2").Value = 3333.3333
$ws.Range("M2").Value = -745.1429000000001
$ws.Range("N2").Value = -3557.3333

$ws.Range("H22").Value = 570.8333
$ws.Range("I22").Value = 391.66666
$ws.Range("K22").Value = 391.66666
$ws.Range("M22").Value = -96.66665999999998

$ws.Range("H27").Value = 570.8333
$ws.Range("I27").Value = 391.66666
$ws.Range("K27").Value = 391.66666
$ws.Range("M27").Value = -284.66666

$ws.Range("H46").Value = 5973.718
$ws.Range("J46").Value = 6131.5
$ws.Range("L46").Value = 6131.5
$ws.Range("N46").Value = -6507.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 81040.2
$ws.Range("J80").Value = 81040.2
$ws.Range("L80").Value = 81040.2
$ws.Range("N80").Value = -83036.2

$ws.Range("H83").Value = 81040.2
$ws.Range("J83").Value = 81040.2
$ws.Range("L83").Value = 243120.6
$ws.Range("N83").Value = -253104.6

$ws.Range("H136").Value = 6814.774
$ws.Range("I136").Value = 1733.5
$ws.Range("J136").Value = 9941.712
$ws.Range("K136").Value = 5200.5
$ws.Range("L136").Value = 29825.136
$ws.Range("M136").Value = -2650.5
$ws.Range("N136").Value = -34925.136
